$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restyle the status legend box (F2:G5) with a thin/medium box border,
#     Arial 10pt font, and wrap text on the "status" column (F) ---
$rngF = $ws.Range("F2:F5")
$rngF.Font.Name = "Arial"
$rngF.Font.Size = 10
$rngF.Font.ColorIndex = 8
$rngF.Borders.LineStyle = 1
$rngF.Borders.Weight = -4138
$rngF.WrapText = $true

$rngG = $ws.Range("G2:G5")
$rngG.Font.Name = "Arial"
$rngG.Font.Size = 10
$rngG.Font.ColorIndex = 8
$rngG.Borders.LineStyle = 1
$rngG.Borders.Weight = -4138

# --- Clear the stray leftover legend text in F6:G7 (the empty-legend bug
#     this commit fixes), restoring plain default formatting ---
$ws.Range("A1").Copy()
$ws.Range("F6:G7").PasteSpecial(-4122)
$ws.Range("F6:G7").ClearContents()
$excel.CutCopyMode = $false

# --- Add new task row 11 ("test" / "Be future") ---
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A11").Value2 = "test"
$ws.Range("C11").Value2 = "Be future"
$ws.Rows.Item(11).RowHeight = 15.75

# --- Selection moves to E10 ---
$null = $ws.Range("E10").Select()

# --- Page setup (portrait, paper size 285) ---
$ws.PageSetup.PaperSize = 285
$ws.PageSetup.Orientation = 1
